# Generate Report for Handback
# This script updates the handback-status workbook so that the
# "109dcc50-..." / "a3dae023-..." handoff identifiers (and their
# associated hash/timestamp values) are replaced by the new
# "eb1d9865-..." / "ffff5c3722a8-..." identifiers, consolidating the
# per-row hash/timestamp values so both rows share the same values.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")

$url2 = "https://github.com/OpenLocalizationTest/oltest/blob/332d22984f586bb1cf73b6fb171475f982528f71/e2e/109dcc50-6bd8-4e8b-b586-bcfafe5be738.md"
$url3 = "https://github.com/OpenLocalizationTest/oltest/blob/332d22984f586bb1cf73b6fb171475f982528f71/e2e/a3dae023-59cc-481c-b11e-951bd45f84ca.md"

$ws.Hyperlinks.Delete()

$ws.Range("A2").Value = "eb1d9865-352f-4f8d-8df1-db1d8c63ef0b.md"
$ws.Range("A3").Value = "ffff5c3722a8-9abb-4a18-9c82-bf8ebf34e2d4.md"

$ws.Hyperlinks.Add($ws.Range("A2"), $url2, "", "", "eb1d9865-352f-4f8d-8df1-db1d8c63ef0b.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A3"), $url3, "", "", "ffff5c3722a8-9abb-4a18-9c82-bf8ebf34e2d4.md") | Out-Null

# ---------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

$urlA2 = "https://github.com/OpenLocalizationTest/oltest/blob/332d22984f586bb1cf73b6fb171475f982528f71/e2e/109dcc50-6bd8-4e8b-b586-bcfafe5be738.md"
$urlB2 = "https://github.com/OpenLocalizationTest/oltest/blob/332d22984f586bb1cf73b6fb171475f982528f71/e2e/109dcc50-6bd8-4e8b-b586-bcfafe5be738.md"
$urlD2 = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/4a53dabaacb2427959445960a82da41bd3a07e78/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/109dcc50-6bd8-4e8b-b586-bcfafe5be738.15ebeab342c42b42b91b326c5aeb3f11342164bb.zh-cn.xlf"
$urlF2 = "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/e3e4dea91ea0ab397e8f84b9a478388688bee1e0/e2e/109dcc50-6bd8-4e8b-b586-bcfafe5be738.md"
$urlG2 = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/fa3fbcc540fa84c76a2336bb8d41fd762552e39e/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/109dcc50-6bd8-4e8b-b586-bcfafe5be738.15ebeab342c42b42b91b326c5aeb3f11342164bb.zh-cn.xlf"
$urlA3 = "https://github.com/OpenLocalizationTest/oltest/blob/332d22984f586bb1cf73b6fb171475f982528f71/e2e/a3dae023-59cc-481c-b11e-951bd45f84ca.md"
$urlB3 = "https://github.com/OpenLocalizationTest/oltest/blob/332d22984f586bb1cf73b6fb171475f982528f71/e2e/a3dae023-59cc-481c-b11e-951bd45f84ca.md"
$urlD3 = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/4a53dabaacb2427959445960a82da41bd3a07e78/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/a3dae023-59cc-481c-b11e-951bd45f84ca.b46bf5b7f4eb4de4fdd58708686169a2bef7598a.zh-cn.xlf"
$urlF3 = "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/e3e4dea91ea0ab397e8f84b9a478388688bee1e0/e2e/a3dae023-59cc-481c-b11e-951bd45f84ca.md"
$urlG3 = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/fa3fbcc540fa84c76a2336bb8d41fd762552e39e/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/a3dae023-59cc-481c-b11e-951bd45f84ca.b46bf5b7f4eb4de4fdd58708686169a2bef7598a.zh-cn.xlf"

$ws.Hyperlinks.Delete()

$newMd1 = "eb1d9865-352f-4f8d-8df1-db1d8c63ef0b.md"
$newMd2 = "ffff5c3722a8-9abb-4a18-9c82-bf8ebf34e2d4.md"
$newXlfZh = "eb1d9865-352f-4f8d-8df1-db1d8c63ef0b.a9aebbb1fc10ccba2e3bc55dd1fbb73fd32cf529.zh-cn.xlf"
$handoffDtZh = "2016-03-20 08:50:39"
$handbackDtZh = "2016-03-20 08:50:58"

$ws.Range("A2").Value = $newMd1
$ws.Range("B2").Value = ".md"
$ws.Range("C2").Value = "Handed back: in sync with en-US"
$ws.Range("D2").Value = $newXlfZh
$ws.Range("E2").Value = $handoffDtZh
$ws.Range("F2").Value = $newMd1
$ws.Range("G2").Value = $newXlfZh
$ws.Range("H2").Value = $handbackDtZh
$ws.Range("I2").Value = "Include"

$ws.Range("A3").Value = $newMd2
$ws.Range("B3").Value = ".md"
$ws.Range("C3").Value = "Handed back: in sync with en-US"
$ws.Range("D3").Value = $newXlfZh
$ws.Range("E3").Value = $handoffDtZh
$ws.Range("F3").Value = $newMd2
$ws.Range("G3").Value = $newXlfZh
$ws.Range("H3").Value = $handbackDtZh
$ws.Range("I3").Value = "Include"

$ws.Hyperlinks.Add($ws.Range("A2"), $urlA2, "", "", $newMd1) | Out-Null
$ws.Hyperlinks.Add($ws.Range("B2"), $urlB2, "", "", ".md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D2"), $urlD2, "", "", $newXlfZh) | Out-Null
$ws.Hyperlinks.Add($ws.Range("F2"), $urlF2, "", "", $newMd1) | Out-Null
$ws.Hyperlinks.Add($ws.Range("G2"), $urlG2, "", "", $newXlfZh) | Out-Null
$ws.Hyperlinks.Add($ws.Range("A3"), $urlA3, "", "", $newMd2) | Out-Null
$ws.Hyperlinks.Add($ws.Range("B3"), $urlB3, "", "", ".md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D3"), $urlD3, "", "", $newXlfZh) | Out-Null
$ws.Hyperlinks.Add($ws.Range("F3"), $urlF3, "", "", $newMd2) | Out-Null
$ws.Hyperlinks.Add($ws.Range("G3"), $urlG3, "", "", $newXlfZh) | Out-Null

# ---------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")

$urlA2 = "https://github.com/OpenLocalizationTest/oltest/blob/332d22984f586bb1cf73b6fb171475f982528f71/e2e/109dcc50-6bd8-4e8b-b586-bcfafe5be738.md"
$urlB2 = "https://github.com/OpenLocalizationTest/oltest/blob/332d22984f586bb1cf73b6fb171475f982528f71/e2e/109dcc50-6bd8-4e8b-b586-bcfafe5be738.md"
$urlD2 = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/64650fa90b9746eaa4376300910808b45395f0de/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/109dcc50-6bd8-4e8b-b586-bcfafe5be738.15ebeab342c42b42b91b326c5aeb3f11342164bb.de-de.xlf"
$urlF2 = "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/0210b075af18aca75a5f5bb6c953499c6a2d970b/e2e/109dcc50-6bd8-4e8b-b586-bcfafe5be738.md"
$urlG2 = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/f079772503cbc40b433ebdf88d1353ee81a4840c/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/109dcc50-6bd8-4e8b-b586-bcfafe5be738.15ebeab342c42b42b91b326c5aeb3f11342164bb.de-de.xlf"
$urlA3 = "https://github.com/OpenLocalizationTest/oltest/blob/332d22984f586bb1cf73b6fb171475f982528f71/e2e/a3dae023-59cc-481c-b11e-951bd45f84ca.md"
$urlB3 = "https://github.com/OpenLocalizationTest/oltest/blob/332d22984f586bb1cf73b6fb171475f982528f71/e2e/a3dae023-59cc-481c-b11e-951bd45f84ca.md"
$urlD3 = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/64650fa90b9746eaa4376300910808b45395f0de/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/a3dae023-59cc-481c-b11e-951bd45f84ca.b46bf5b7f4eb4de4fdd58708686169a2bef7598a.de-de.xlf"
$urlF3 = "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/0210b075af18aca75a5f5bb6c953499c6a2d970b/e2e/a3dae023-59cc-481c-b11e-951bd45f84ca.md"
$urlG3 = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/f079772503cbc40b433ebdf88d1353ee81a4840c/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/a3dae023-59cc-481c-b11e-951bd45f84ca.b46bf5b7f4eb4de4fdd58708686169a2bef7598a.de-de.xlf"

$ws.Hyperlinks.Delete()

$newXlfDe = "eb1d9865-352f-4f8d-8df1-db1d8c63ef0b.a9aebbb1fc10ccba2e3bc55dd1fbb73fd32cf529.de-de.xlf"
$handoffDtDe = "2016-03-20 08:50:42"
$handbackDtDe = "2016-03-20 08:51:05"

$ws.Range("A2").Value = $newMd1
$ws.Range("B2").Value = ".md"
$ws.Range("C2").Value = "Handed back: in sync with en-US"
$ws.Range("D2").Value = $newXlfDe
$ws.Range("E2").Value = $handoffDtDe
$ws.Range("F2").Value = $newMd1
$ws.Range("G2").Value = $newXlfDe
$ws.Range("H2").Value = $handbackDtDe
$ws.Range("I2").Value = "Include"

$ws.Range("A3").Value = $newMd2
$ws.Range("B3").Value = ".md"
$ws.Range("C3").Value = "Handed back: in sync with en-US"
$ws.Range("D3").Value = $newXlfDe
$ws.Range("E3").Value = $handoffDtDe
$ws.Range("F3").Value = $newMd2
$ws.Range("G3").Value = $newXlfDe
$ws.Range("H3").Value = $handbackDtDe
$ws.Range("I3").Value = "Include"

$ws.Hyperlinks.Add($ws.Range("A2"), $urlA2, "", "", $newMd1) | Out-Null
$ws.Hyperlinks.Add($ws.Range("B2"), $urlB2, "", "", ".md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D2"), $urlD2, "", "", $newXlfDe) | Out-Null
$ws.Hyperlinks.Add($ws.Range("F2"), $urlF2, "", "", $newMd1) | Out-Null
$ws.Hyperlinks.Add($ws.Range("G2"), $urlG2, "", "", $newXlfDe) | Out-Null
$ws.Hyperlinks.Add($ws.Range("A3"), $urlA3, "", "", $newMd2) | Out-Null
$ws.Hyperlinks.Add($ws.Range("B3"), $urlB3, "", "", ".md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D3"), $urlD3, "", "", $newXlfDe) | Out-Null
$ws.Hyperlinks.Add($ws.Range("F3"), $urlF3, "", "", $newMd2) | Out-Null
$ws.Hyperlinks.Add($ws.Range("G3"), $urlG3, "", "", $newXlfDe) | Out-Null

$wb.Save()
